$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inventory")

# Preserve the current last row (row 23: /build.gradle) before the new rows push it down to row 26.
$a23 = $ws.Range("A23").Value2
$f23 = $ws.Range("F23").Value2
$g23 = $ws.Range("G23").Value2

$ws.Range("A26").Value = $a23
$ws.Range("F26").Value = $f23
$ws.Range("G26").Value = $g23

# The old row 23 only keeps B23/G23 populated with new content; A23/F23 become blank.
$ws.Range("A23").ClearContents()
$ws.Range("F23").ClearContents()

# Add the new "/src/test" inventory rows (23-25).
$ws.Range("B23").Value = "[/test](/src/test)"

$ws.Range("C24").Value = "[/ml-config](/src/test/ml-config)"
$ws.Range("F24").Value = "Configuration required to execute unit tests, including roles and users."

$ws.Range("G23").Value = "Conditionally deployed by the ``restrictUnitTestingDeployment`` Gradle task."

$ws.Range("C25").Value = "[/ml-modules](/src/test/ml-modules)"
$ws.Range("F25").Value = "All test suites, which can also include test data."
$ws.Range("G25").Value = "[MarkLogic Unit Test user guide](https://marklogic-community.github.io/marklogic-unit-test/)"
